$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.686.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.443.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +3.08%  "

$ws.Range("E10").Value = "  -2.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.407"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.028.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.37%  "

$ws.Range("E13").Value = "  +2.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.447.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("E16").Value = "  -1.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.730.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("E18").Value = "  +1.68%  "

$ws.Range("E19").Value = "  +1.50%  "

$ws.Range("E20").Value = "  -3.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.11%  "

$ws.Range("E22").Value = "  -0.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("E25").Value = "  -2.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.589.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.182"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  -2.96%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  -8.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.54%  "

$ws.Range("E35").Value = "  +2.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "31.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.57%  "

$ws.Range("E38").Value = "  -2.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "170.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.476.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0776"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.789"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.19%  "

$ws.Range("E45").Value = "  -3.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.561.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.29%  "

$ws.Range("E51").Value = "  -0.09%  "
